$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Probabilities")

$ws.Cells.Item(2, 2).Value = "2025-11-30T10:00:00"
$ws.Cells.Item(2, 3).Value = "Адмирал"
$ws.Cells.Item(2, 4).Value = "Динамо Мн"
$ws.Cells.Item(2, 5).Value = 897831
$ws.Cells.Item(2, 6).Value = "https://text.khl.ru/text/897831.html"
$ws.Cells.Item(2, 7).Value = 2.617412
$ws.Cells.Item(2, 8).Value = 4.5
$ws.Cells.Item(2, 9).Value = 1.125
$ws.Cells.Item(2, 10).Value = 1.540586
$ws.Cells.Item(2, 11).Value = 2.078999
$ws.Cells.Item(2, 12).Value = 2.8125
$ws.Cells.Item(2, 13).Value = 7.117412
$ws.Cells.Item(2, 14).Value = 28.706351
$ws.Cells.Item(2, 15).Value = 36.917927
$ws.Cells.Item(2, 16).Value = 65.624278
$ws.Cells.Item(2, 17).Value = -0.041041
$ws.Cells.Item(2, 18).Value = 0.2
$ws.Cells.Item(2, 19).Value = 0.284707
$ws.Cells.Item(2, 20).Value = 0.176781
$ws.Cells.Item(2, 21).Value = 0.53833
$ws.Cells.Item(2, 22).Value = 0.280588
$ws.Cells.Item(2, 23).Value = 0.71923
$ws.Cells.Item(2, 24).Value = 0.459734
$ws.Cells.Item(2, 25).Value = 0.540084
$ws.Cells.Item(2, 26).Value = 0.634992
$ws.Cells.Item(2, 27).Value = 0.364827
$ws.Cells.Item(2, 28).Value = 0.77787
$ws.Cells.Item(2, 29).Value = 0.221948
$ws.Cells.Item(2, 30).Value = 0.877712
$ws.Cells.Item(2, 31).Value = 0.122106
$ws.Cells.Item(2, 32).Value = 0.614955
$ws.Cells.Item(2, 33).Value = 0.385045
$ws.Cells.Item(2, 34).Value = 0.344695
$ws.Cells.Item(2, 35).Value = 0.655305
$ws.Cells.Item(2, 36).Value = 0.771042
$ws.Cells.Item(2, 37).Value = 0.228958
$ws.Cells.Item(2, 38).Value = 0.533521
$ws.Cells.Item(2, 39).Value = 0.466479
$ws.Cells.Item(2, 40).Value = 0.644371
$ws.Cells.Item(2, 41).Value = 0.850298

$ws.Cells.Item(3, 2).Value = "2025-11-30T10:00:00"
$ws.Cells.Item(3, 3).Value = "Амур"
$ws.Cells.Item(3, 4).Value = "ХК Сочи"
$ws.Cells.Item(3, 5).Value = 897832
$ws.Cells.Item(3, 6).Value = "https://text.khl.ru/text/897832.html"
$ws.Cells.Item(3, 7).Value = 1.225758
$ws.Cells.Item(3, 8).Value = 0.961538
$ws.Cells.Item(3, 9).Value = 1.868405
$ws.Cells.Item(3, 10).Value = 7.038462
$ws.Cells.Item(3, 11).Value = 4.13211
$ws.Cells.Item(3, 12).Value = 1.414972
$ws.Cells.Item(3, 13).Value = 2.187297
$ws.Cells.Item(3, 14).Value = 24.98031
$ws.Cells.Item(3, 15).Value = 25.615145
$ws.Cells.Item(3, 16).Value = 50.595455
$ws.Cells.Item(3, 17).Value = -0.191993
$ws.Cells.Item(3, 18).Value = -0.2
$ws.Cells.Item(3, 19).Value = 0.827776
$ws.Cells.Item(3, 20).Value = 0.091762
$ws.Cells.Item(3, 21).Value = 0.07685
$ws.Cells.Item(3, 22).Value = 0.196421
$ws.Cells.Item(3, 23).Value = 0.799967
$ws.Cells.Item(3, 24).Value = 0.350229
$ws.Cells.Item(3, 25).Value = 0.646159
$ws.Cells.Item(3, 26).Value = 0.520867
$ws.Cells.Item(3, 27).Value = 0.475522
$ws.Cells.Item(3, 28).Value = 0.678623
$ws.Cells.Item(3, 29).Value = 0.317765
$ws.Cells.Item(3, 30).Value = 0.803636
$ws.Cells.Item(3, 31).Value = 0.192752
$ws.Cells.Item(3, 32).Value = 0.917635
$ws.Cells.Item(3, 33).Value = 0.082365
$ws.Cells.Item(3, 34).Value = 0.780622
$ws.Cells.Item(3, 35).Value = 0.219378
$ws.Cells.Item(3, 36).Value = 0.413325
$ws.Cells.Item(3, 37).Value = 0.586675
$ws.Cells.Item(3, 38).Value = 0.170132
$ws.Cells.Item(3, 39).Value = 0.829868
$ws.Cells.Item(3, 40).Value = 0.967298
$ws.Cells.Item(3, 41).Value = 0.308085

$ws.Cells.Item(4, 2).Value = "2025-11-30T17:00:00"
$ws.Cells.Item(4, 3).Value = "Ак Барс"
$ws.Cells.Item(4, 4).Value = "Драконы"
$ws.Cells.Item(4, 5).Value = 897833
$ws.Cells.Item(4, 6).Value = "https://text.khl.ru/text/897833.html"
$ws.Cells.Item(4, 7).Value = 3.055625
$ws.Cells.Item(4, 8).Value = 3.676454
$ws.Cells.Item(4, 9).Value = 3.88627
$ws.Cells.Item(4, 10).Value = 5.5
$ws.Cells.Item(4, 11).Value = 4.277812
$ws.Cells.Item(4, 12).Value = 3.781362
$ws.Cells.Item(4, 13).Value = 6.732079
$ws.Cells.Item(4, 14).Value = 33.778683
$ws.Cells.Item(4, 15).Value = 31.087527
$ws.Cells.Item(4, 16).Value = 64.866211
$ws.Cells.Item(4, 17).Value = -0.007958
$ws.Cells.Item(4, 18).Value = 0.151468
$ws.Cells.Item(4, 19).Value = 0.49328
$ws.Cells.Item(4, 20).Value = 0.140849
$ws.Cells.Item(4, 21).Value = 0.359373
$ws.Cells.Item(4, 22).Value = 0.040717
$ws.Cells.Item(4, 23).Value = 0.952785
$ws.Cells.Item(4, 24).Value = 0.096294
$ws.Cells.Item(4, 25).Value = 0.897208
$ws.Cells.Item(4, 26).Value = 0.185875
$ws.Cells.Item(4, 27).Value = 0.807627
$ws.Cells.Item(4, 28).Value = 0.3062
$ws.Cells.Item(4, 29).Value = 0.687302
$ws.Cells.Item(4, 30).Value = 0.444732
$ws.Cells.Item(4, 31).Value = 0.548771
$ws.Cells.Item(4, 32).Value = 0.926781
$ws.Cells.Item(4, 33).Value = 0.073219
$ws.Cells.Item(4, 34).Value = 0.799846
$ws.Cells.Item(4, 35).Value = 0.200154
$ws.Cells.Item(4, 36).Value = 0.891025
$ws.Cells.Item(4, 37).Value = 0.108975
$ws.Cells.Item(4, 38).Value = 0.72808
$ws.Cells.Item(4, 39).Value = 0.27192
$ws.Cells.Item(4, 40).Value = 0.758018
$ws.Cells.Item(4, 41).Value = 0.640376
